$d = $word.ActiveDocument

# 1. Merge "Introduction" + " To Build Responsive" into one run's text.
#    (Find/Replace on the visible text achieves the same resulting text content.)
$d.Content.Find.Execute("Introduction To Build Responsive", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Introduction To Build Responsive", 0)

# 2. Replace "Loso" with "this " (note trailing space) in the icons paragraph.
$d.Content.Find.Execute("Loso", $true, $false, $false, $false, $false,
                         $true, 1, $false, "this ", 2)
